$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ece score
$ws.Range("C4").Value = 2.844
$ws.Range("D4").Value = 3.105
$ws.Range("E4").Value = 2.388
$ws.Range("F4").Value = 2.464
$ws.Range("G4").Value = 2.141
$ws.Range("H4").Value = 3.269

$ws.Range("C5").Value = 2.672
$ws.Range("D5").Value = 3.425
$ws.Range("E5").Value = 3.443
$ws.Range("F5").Value = 2.414
$ws.Range("G5").Value = 3.251
$ws.Range("H5").Value = 3.201

$ws.Range("C6").Value = 3.008
$ws.Range("D6").Value = 3.318
$ws.Range("E6").Value = 3.518
$ws.Range("F6").Value = 2.818
$ws.Range("G6").Value = 3.101
$ws.Range("H6").Value = 3.262

# brier score loss
$ws.Range("C7").Value = 0.752
$ws.Range("D7").Value = 0.461
$ws.Range("E7").Value = 0.443
$ws.Range("F7").Value = 0.648
$ws.Range("G7").Value = 0.488
$ws.Range("H7").Value = 0.428

$ws.Range("C8").Value = 0.801
$ws.Range("D8").Value = 0.541
$ws.Range("E8").Value = 0.503
$ws.Range("F8").Value = 0.763
$ws.Range("G8").Value = 0.6
$ws.Range("H8").Value = 0.524

$ws.Range("C9").Value = 0.618
$ws.Range("D9").Value = 0.477
$ws.Range("E9").Value = 0.439
$ws.Range("F9").Value = 0.677
$ws.Range("G9").Value = 0.542
$ws.Range("H9").Value = 0.483
